$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5955141.5
$ws.Range("I33").Value = 1802489.1
$ws.Range("J33").Value = 20835480
$ws.Range("K33").Value = 1802489.1
$ws.Range("L33").Value = 20835480
$ws.Range("M33").Value = -1802260.1
$ws.Range("N33").Value = -20835938
$ws.Range("H80").Value = 613.05884
$ws.Range("I80").Value = 443
$ws.Range("J80").Value = 856
$ws.Range("K80").Value = 1329
$ws.Range("L80").Value = 2568
$ws.Range("M80").Value = -331
$ws.Range("N80").Value = -4564
$ws.Range("H83").Value = 613.05884
$ws.Range("I83").Value = 443
$ws.Range("J83").Value = 856
$ws.Range("K83").Value = 3987
$ws.Range("L83").Value = 7704
$ws.Range("M83").Value = 1005
$ws.Range("N83").Value = -17688
$ws.Range("H112").Value = 2665.3572
$ws.Range("J112").Value = 2642.5908
$ws.Range("L112").Value = 7927.7724
$ws.Range("N112").Value = -10143.7724
$ws.Range("H113").Value = 3901.25
$ws.Range("I113").Value = 3535
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 3535
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -281
$ws.Range("N113").Value = -11508
$ws.Range("H131").Value = 4385.174
$ws.Range("I131").Value = 3857.2273
$ws.Range("J131").Value = 16000
$ws.Range("K131").Value = 11571.6819
$ws.Range("L131").Value = 48000
$ws.Range("M131").Value = -6531.6819
$ws.Range("N131").Value = -58080
$ws.Range("H133").Value = 199998.67
$ws.Range("J133").Value = 199998.4
$ws.Range("L133").Value = 199998.4
$ws.Range("N133").Value = -210118.4
$ws.Range("H137").Value = 2280112.5
$ws.Range("I137").Value = 7401.136
$ws.Range("J137").Value = 4552824
$ws.Range("K137").Value = 22203.408
$ws.Range("L137").Value = 13658472
$ws.Range("M137").Value = -19653.408
$ws.Range("N137").Value = -13663572
$ws.Range("H138").Value = 4323.603
$ws.Range("I138").Value = 4258.3887
$ws.Range("J138").Value = 4396.9688
$ws.Range("K138").Value = 12775.1661
$ws.Range("L138").Value = 13190.9064
$ws.Range("M138").Value = -7635.166100000002
$ws.Range("N138").Value = -23470.9064

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 76216.5
$ws.Range("I45").Value = 87961
$ws.Range("K45").Value = 87961
$ws.Range("M45").Value = -87584
$ws.Range("H61").Value = 2398530.5
$ws.Range("I61").Value = 72851.13
$ws.Range("J61").Value = 5887049.5
$ws.Range("K61").Value = 72851.13
$ws.Range("L61").Value = 5887049.5
$ws.Range("M61").Value = -72639.13
$ws.Range("N61").Value = -5887473.5
$ws.Range("H122").Value = 3882.0208
$ws.Range("I122").Value = 3188.158
$ws.Range("J122").Value = 4336.6206
$ws.Range("K122").Value = 9564.474
$ws.Range("L122").Value = 13009.8618
$ws.Range("M122").Value = -7114.474
$ws.Range("N122").Value = -17909.8618
$ws.Range("H132").Value = 1568.8096
$ws.Range("I132").Value = 1172.8379
$ws.Range("K132").Value = 3518.5137
$ws.Range("M132").Value = -988.5137
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 2398530.5
$ws.Range("I136").Value = 72851.13
$ws.Range("J136").Value = 5887049.5
$ws.Range("K136").Value = 218553.39
$ws.Range("L136").Value = 17661148.5
$ws.Range("M136").Value = -216003.39
$ws.Range("N136").Value = -17666248.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2696.9014
$ws.Range("I31").Value = 2117.4375
$ws.Range("J31").Value = 3172.359
$ws.Range("K31").Value = 2117.4375
$ws.Range("L31").Value = 3172.359
$ws.Range("M31").Value = -1822.4375
$ws.Range("N31").Value = -3762.359
$ws.Range("H34").Value = 2696.9014
$ws.Range("I34").Value = 2117.4375
$ws.Range("J34").Value = 3172.359
$ws.Range("K34").Value = 2117.4375
$ws.Range("L34").Value = 3172.359
$ws.Range("M34").Value = -1915.4375
$ws.Range("N34").Value = -3576.359
$ws.Range("H99").Value = 62502748
$ws.Range("I99").Value = 2195.2
$ws.Range("J99").Value = 166670340
$ws.Range("K99").Value = 2195.2
$ws.Range("L99").Value = 166670340
$ws.Range("M99").Value = -697.1999999999998
$ws.Range("N99").Value = -166673336
$ws.Range("H122").Value = 1722.9166
$ws.Range("I122").Value = 1674.8462
$ws.Range("J122").Value = 1847.9
$ws.Range("K122").Value = 5024.5386
$ws.Range("L122").Value = 5543.700000000001
$ws.Range("M122").Value = -2574.5386
$ws.Range("N122").Value = -10443.7
$ws.Range("H126").Value = 62502748
$ws.Range("I126").Value = 2195.2
$ws.Range("J126").Value = 166670340
$ws.Range("K126").Value = 6585.599999999999
$ws.Range("L126").Value = 500011020
$ws.Range("M126").Value = -4115.599999999999
$ws.Range("N126").Value = -500015960
$ws.Range("H134").Value = 3233.2307
$ws.Range("I134").Value = 3086
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 9258
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -6723
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1984490.1
$ws.Range("J122").Value = 7936634
$ws.Range("L122").Value = 71429706
$ws.Range("N122").Value = -71434606

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4483.7144
$ws.Range("I70").Value = 4703.3335
$ws.Range("J70").Value = 3166
$ws.Range("K70").Value = 4703.3335
$ws.Range("L70").Value = 3166
$ws.Range("M70").Value = -4433.3335
$ws.Range("N70").Value = -3706
$ws.Range("H73").Value = 4483.7144
$ws.Range("I73").Value = 4703.3335
$ws.Range("J73").Value = 3166
$ws.Range("K73").Value = 4703.3335
$ws.Range("L73").Value = 3166
$ws.Range("M73").Value = -3767.3335
$ws.Range("N73").Value = -5038
$ws.Range("H122").Value = 2444.6843
$ws.Range("I122").Value = 2290.5625
$ws.Range("K122").Value = 6871.6875
$ws.Range("M122").Value = -4421.6875
$ws.Range("H123").Value = 97999.5
$ws.Range("J123").Value = 97999.5
$ws.Range("L123").Value = 97999.5
$ws.Range("N123").Value = -102899.5
$ws.Range("H126").Value = 2650.3635
$ws.Range("I126").Value = 2098.6924
$ws.Range("J126").Value = 3447.2222
$ws.Range("K126").Value = 6296.0772
$ws.Range("L126").Value = 10341.6666
$ws.Range("M126").Value = -3826.0772
$ws.Range("N126").Value = -15281.6666
$ws.Range("H132").Value = 14511414
$ws.Range("I132").Value = 3985.2222
$ws.Range("J132").Value = 27568100
$ws.Range("K132").Value = 11955.6666
$ws.Range("L132").Value = 82704300
$ws.Range("M132").Value = -9425.6666
$ws.Range("N132").Value = -82709360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3305.9092
$ws.Range("I122").Value = 3052
$ws.Range("K122").Value = 9156
$ws.Range("M122").Value = -6706
$ws.Range("H132").Value = 4307.1313
$ws.Range("I132").Value = 2756.682
$ws.Range("K132").Value = 8270.045999999998
$ws.Range("M132").Value = -5740.045999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H81").Value = 37854.605
$ws.Range("I81").Value = 2363.2
$ws.Range("K81").Value = 4726.4
$ws.Range("M81").Value = -3665.4
$ws.Range("H84").Value = 37854.605
$ws.Range("I84").Value = 2363.2
$ws.Range("K84").Value = 23632
$ws.Range("M84").Value = -18328
$ws.Range("H132").Value = 1838.5454
$ws.Range("I132").Value = 1528.625
$ws.Range("J132").Value = 2665
$ws.Range("K132").Value = 4585.875
$ws.Range("L132").Value = 7995
$ws.Range("M132").Value = -2055.875
$ws.Range("N132").Value = -13055

